$wb = $excel.ActiveWorkbook

# ---- Sheet: Metadata ----
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("A2").Value = "29 Oct 2025, 06:31 PM"

# ---- Sheet: Top Gainers ----
$ws = $wb.Worksheets.Item("Top Gainers")
$ws.Range("C57").Value = 3.9981
$ws.Range("D57").Value = 7.4592
$ws.Range("E57").Value = 27.1054
$ws.Range("B61").Value = "APARINDS"
$ws.Range("C61").Value = 3.8924
$ws.Range("D61").Value = 8.3414
$ws.Range("E61").Value = 15.5876
$ws.Range("B62").Value = "HITECHGEAR"
$ws.Range("C62").Value = 3.8587
$ws.Range("D62").Value = 1.1486
$ws.Range("E62").Value = 9.9254
$ws.Range("B63").Value = "ORIENTTECH"
$ws.Range("C63").Value = 3.827
$ws.Range("D63").Value = 0.5247000000000001
$ws.Range("E63").Value = 32.6784
$ws.Range("B64").Value = "ICRA"
$ws.Range("C64").Value = 3.7985
$ws.Range("D64").Value = 4.4793
$ws.Range("E64").Value = 2.8828
$ws.Range("B65").Value = "SALASAR"
$ws.Range("C65").Value = 3.7935
$ws.Range("D65").Value = 4.7872
$ws.Range("E65").Value = 11.0485
$ws.Range("B66").Value = "NPST"
$ws.Range("C66").Value = 3.7841
$ws.Range("D66").Value = -2.0689
$ws.Range("E66").Value = -3.5677
$ws.Range("B67").Value = "DCW"
$ws.Range("C67").Value = 3.7544
$ws.Range("D67").Value = 2.3219
$ws.Range("E67").Value = -3.9753
$ws.Range("B68").Value = "RHETAN"
$ws.Range("C68").Value = 3.754
$ws.Range("D68").Value = 4.178
$ws.Range("E68").Value = 6.549
$ws.Range("B69").Value = "HINDPETRO"
$ws.Range("C69").Value = 3.6935
$ws.Range("D69").Value = 6.9335
$ws.Range("E69").Value = 5.7397
$ws.Range("B70").Value = "BHARTIHEXA"
$ws.Range("C70").Value = 3.6718
$ws.Range("D70").Value = 7.0877
$ws.Range("E70").Value = 15.3332
$ws.Range("B71").Value = "HLEGLAS"
$ws.Range("C71").Value = 3.659
$ws.Range("D71").Value = 8.115500000000001
$ws.Range("E71").Value = 27.1239
$ws.Range("B72").Value = "RHIM"
$ws.Range("C72").Value = 3.6544
$ws.Range("D72").Value = 3.2276
$ws.Range("E72").Value = 5.1826
$ws.Range("B73").Value = "SHK"
$ws.Range("C73").Value = 3.6347
$ws.Range("D73").Value = 2.388
$ws.Range("E73").Value = -1.932
$ws.Range("B74").Value = "BCLIND"
$ws.Range("C74").Value = 3.6271
$ws.Range("D74").Value = 2.2945
$ws.Range("E74").Value = 0.1728
$ws.Range("B75").Value = "MUKANDLTD"
$ws.Range("C75").Value = 3.6133
$ws.Range("D75").Value = 11.9685
$ws.Range("E75").Value = 9.550800000000001
$ws.Range("B76").Value = "CGPOWER"
$ws.Range("C76").Value = 3.6125
$ws.Range("D76").Value = 3.4192
$ws.Range("E76").Value = 1.0325

# ---- Sheet: 1 Month Performance ----
$ws = $wb.Worksheets.Item("1 Month Performance")
$ws.Range("C3").Value = 79.9766
$ws.Range("B5").Value = "MCLEODRUSS"
$ws.Range("C5").Value = 69.8655
$ws.Range("B6").Value = "PROZONER"
$ws.Range("C6").Value = 68.3711
$ws.Range("B7").Value = "IFBAGRO"
$ws.Range("C7").Value = 66.1892
$ws.Range("B8").Value = "BGRENERGY"
$ws.Range("C8").Value = 64.929
$ws.Range("B9").Value = "ESSARSHPNG"
$ws.Range("C9").Value = 64.91160000000001
$ws.Range("B10").Value = "PANACHE"
$ws.Range("C10").Value = 62.6487
$ws.Range("B11").Value = "MAHASTEEL"
$ws.Range("C11").Value = 55.9703
$ws.Range("B12").Value = "INOXGREEN"
$ws.Range("C12").Value = 51.0181
$ws.Range("B13").Value = "STALLION"
$ws.Range("C13").Value = 46.4325
$ws.Range("B14").Value = "ORIENTTECH"
$ws.Range("C14").Value = 45.3321
$ws.Range("B15").Value = "TVSSRICHAK"
$ws.Range("C15").Value = 40.7778
$ws.Range("B16").Value = "MTARTECH"
$ws.Range("C16").Value = 40.7213
$ws.Range("B17").Value = "SEJALLTD"
$ws.Range("C17").Value = 37.4301
$ws.Range("B18").Value = "V2RETAIL"
$ws.Range("C18").Value = 37.2004
$ws.Range("B19").Value = "RAMAPHO"
$ws.Range("C19").Value = 36.9731
$ws.Range("B20").Value = "SANDUMA"
$ws.Range("C20").Value = 36.9057
$ws.Range("B21").Value = "TARACHAND"
$ws.Range("C21").Value = 36.4813
$ws.Range("B22").Value = "NETWEB"
$ws.Range("C22").Value = 36.1199
$ws.Range("B23").Value = "SAMMAANCAP"
$ws.Range("C23").Value = 35.5128
$ws.Range("B24").Value = "ONMOBILE"
$ws.Range("C24").Value = 35.4702
$ws.Range("B25").Value = "SHAREINDIA"
$ws.Range("C25").Value = 35.3207
$ws.Range("B26").Value = "SOUTHBANK"
$ws.Range("C26").Value = 35.2819
$ws.Range("B27").Value = "TVSELECT"
$ws.Range("C27").Value = 35.1983
$ws.Range("B28").Value = "RAMCOSYS"
$ws.Range("C28").Value = 34.6928
$ws.Range("B29").Value = "MAANALU"
$ws.Range("C29").Value = 34.4803
$ws.Range("B30").Value = "MEGASOFT"
$ws.Range("C30").Value = 33.4399
$ws.Range("B31").Value = "BHARATSE"
$ws.Range("C31").Value = 31.8611
$ws.Range("B32").Value = "EMKAY"
$ws.Range("C32").Value = 30.3743
$ws.Range("B33").Value = "ATHERENERG"
$ws.Range("C33").Value = 29.116
$ws.Range("B34").Value = "TATVA"
$ws.Range("C34").Value = 28.7451
$ws.Range("B35").Value = "TERASOFT"
$ws.Range("C35").Value = 28.3093
$ws.Range("B36").Value = "CARTRADE"
$ws.Range("C36").Value = 27.5713
$ws.Range("B37").Value = "ARFIN"
$ws.Range("C37").Value = 27.3801
$ws.Range("B38").Value = "MINDTECK"
$ws.Range("C38").Value = 26.9415
$ws.Range("B39").Value = "BHARATWIRE"
$ws.Range("C39").Value = 26.5276
$ws.Range("B40").Value = "HATSUN"
$ws.Range("C40").Value = 26.492
$ws.Range("B41").Value = "INDORAMA"
$ws.Range("C41").Value = 26.4516
$ws.Range("B42").Value = "IFBIND"
$ws.Range("C42").Value = 26.161
$ws.Range("B43").Value = "ADANIPOWER"
$ws.Range("C43").Value = 25.8247
$ws.Range("B44").Value = "AVALON"
$ws.Range("C44").Value = 25.7352
$ws.Range("B45").Value = "MRPL"
$ws.Range("C45").Value = 25.6265
$ws.Range("B46").Value = "HINDCOPPER"
$ws.Range("C46").Value = 25.3164
$ws.Range("B47").Value = "PRECWIRE"
$ws.Range("C47").Value = 24.679
$ws.Range("B48").Value = "SCI"
$ws.Range("C48").Value = 24.132
$ws.Range("B49").Value = "KICL"
$ws.Range("C49").Value = 24.1119
$ws.Range("B50").Value = "SKYGOLD"
$ws.Range("C50").Value = 24.1079
$ws.Range("B51").Value = "DCBBANK"
$ws.Range("C51").Value = 23.8922
$ws.Range("B52").Value = "AUBANK"
$ws.Range("C52").Value = 23.6964
$ws.Range("B53").Value = "ETHOSLTD"
$ws.Range("C53").Value = 23.1527
$ws.Range("B54").Value = "PVP"
$ws.Range("C54").Value = 22.7524
$ws.Range("B55").Value = "INDIANB"
$ws.Range("C55").Value = 22.6689
$ws.Range("B56").Value = "PRIVISCL"
$ws.Range("C56").Value = 22.3984
$ws.Range("B57").Value = "CPEDU"
$ws.Range("C57").Value = 22.3786
$ws.Range("B58").Value = "LORDSCHLO"
$ws.Range("C58").Value = 22.1791
$ws.Range("B59").Value = "GUJTHEM"
$ws.Range("C59").Value = 22.0704
$ws.Range("B60").Value = "SURYODAY"
$ws.Range("C60").Value = 21.8039
$ws.Range("B61").Value = "TDPOWERSYS"
$ws.Range("C61").Value = 21.7743
$ws.Range("B62").Value = "ORBTEXP"
$ws.Range("C62").Value = 21.6115
$ws.Range("B63").Value = "CEATLTD"
$ws.Range("C63").Value = 20.0239
$ws.Range("B64").Value = "ATL"
$ws.Range("C64").Value = 19.9362
$ws.Range("B65").Value = "GRMOVER"
$ws.Range("C65").Value = 19.7859
$ws.Range("B66").Value = "FEDERALBNK"
$ws.Range("C66").Value = 19.6872
$ws.Range("B67").Value = "SUBROS"
$ws.Range("C67").Value = 19.6508
$ws.Range("B68").Value = "USHAMART"
$ws.Range("C68").Value = 19.6172
$ws.Range("B69").Value = "BANKINDIA"
$ws.Range("C69").Value = 19.3067
$ws.Range("B70").Value = "RBLBANK"
$ws.Range("C70").Value = 19.2556
$ws.Range("B71").Value = "MOLDTECH"
$ws.Range("C71").Value = 19.1891
$ws.Range("B72").Value = "THOMASCOTT"
$ws.Range("C72").Value = 19.1649
$ws.Range("B73").Value = "IIFL"
$ws.Range("C73").Value = 18.9853
$ws.Range("B74").Value = "KARURVYSYA"
$ws.Range("C74").Value = 18.8614
$ws.Range("B75").Value = "LUMAXIND"
$ws.Range("C75").Value = 18.8057
$ws.Range("B76").Value = "REPRO"
$ws.Range("C76").Value = 18.689

# ---- Sheet: distance from Dma50 (append new rows) ----
$ws = $wb.Worksheets.Item("distance from Dma50")
$ws.Range("A31").Value = "📈"
$ws.Range("B31").Value = "NIFTYFINSEREXBNK"
$ws.Range("C31").Value = "N/A"
$ws.Range("A32").Value = "📈"
$ws.Range("B32").Value = "NIFTYMSITTELCM"
$ws.Range("C32").Value = "N/A"
$ws.Range("A33").Value = "📈"
$ws.Range("B33").Value = "NIFTYMSFINSERV"
$ws.Range("C33").Value = "N/A"
